$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 18, shifting the existing "season_title" block
# (and everything after it) down by one row.
$ws.Rows.Item(18).Insert()

# Populate the newly inserted row with the new "cancel" / "CANCEL" localization
# key/value pair.
$ws.Range("A18").Value = "cancel"
$ws.Range("B18").Value = "CANCEL"

# Match the formatting used by the rest of column B (wrap text style).
$ws.Range("B18").WrapText = $true

# Update the active selection to the newly inserted row.
$ws.Range("A18").Select()
